# The source diff only changes two kinds of opaque, auto-generated OOXML
# identifiers that Word's object model has never exposed for direct
# assignment:
#   1. w:rsidR="..." stamped on the REF-field runs (w:fldChar/w:instrText).
#   2. The w:id="..." pair on the bookmark1 bookmarkStart/bookmarkEnd.
#
# Both values carry no semantic meaning (no visible text/content changed -
# this is purely a "regenerate the fixture with a newer tool version"
# commit, per the commit message). The Word COM/VBA object model does not
# expose a settable property for either of these (Bookmark has no .ID
# property, and there is no Range/Run property that maps to the raw
# w:rsidR revision-save-id attribute) - this holds in real MS Word just as
# much as here.
#
# The closest faithful, COM-legitimate reproduction of "this bookmark's
# identity was refreshed" is to recreate the bookmark in place: delete it
# and re-add a bookmark with the same name over the same range. That is
# exactly the operation a macro author would use to force Word to mint a
# fresh internal id for it.
$d = $word.ActiveDocument

$bm = $d.Bookmarks("bookmark1")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("bookmark1", $bmRange) | Out-Null
